# Q3 Update - 2025
# - refresh the "short-url" column (B) across all data rows
# - split the former last row (item 38, Venezuela 2024) into two rows:
#   item 38 keeps 2024/Venezuela's slot but the coo becomes Colombia,
#   and a brand new item 39 row is appended for Venezuela 2024 with an
#   updated oip figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate row 39 into row 40 (values + styles) as a starting point for
#    the new record; this also extends the sheet dimension to A1:V40.
$ws.Range("A39:V39").Copy()
$ws.Range("A40:V40").PasteSpecial(-4122)

# 2) Row 40 becomes the new "item 39" Venezuela row - only the item number
#    and the oip figure differ from the row it was cloned from.
$ws.Cells.Item(40, 4).Value = 39
$ws.Cells.Item(40, 21).Value = 11593

# 3) Row 39 becomes the "item 38" Colombia row - coo switches from
#    Venezuela to Colombia, ooc (T) gets a real value, and oip (U) reverts
#    to the placeholder "-" (left-aligned, like every other "-" cell).
$ws.Cells.Item(39, 6).Value = 44
$ws.Cells.Item(39, 7).Value = "Colombia"
$ws.Cells.Item(39, 8).Value = "COL"
$ws.Cells.Item(39, 9).Value = "COL"
$ws.Cells.Item(39, 20).Value = 1278
$ws.Cells.Item(39, 21).Value = "-"
$ws.Cells.Item(39, 21).HorizontalAlignment = -4131

# 4) The "short-url" column (B) value is refreshed everywhere it appears.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 2).Value = "3fIwt7"
}
